$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: remove bold/border/center-alignment style (style index 1) ---
# Clear A1 content entirely (becomes a blank cell) and strip formatting from whole header row
$ws.Range("A1:AG1").ClearFormats()
$ws.Range("A1").Value = ""

# --- Corrected pre/post/total fixation data (rows 3-7) ---
# Row 3: Revisit count
$ws.Range("C3").Value = 14
$ws.Range("D3").Value = 64
$ws.Range("E3").Value = 14
$ws.Range("I3").Value = 16
$ws.Range("K3").Value = 44
$ws.Range("L3").Value = 20
$ws.Range("N3").Value = 10
$ws.Range("T3").Value = 59
$ws.Range("U3").Value = 18
$ws.Range("X3").Value = 26
$ws.Range("AA3").Value = 45
$ws.Range("AB3").Value = 12

# Row 4: Fixation count
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 261
$ws.Range("E4").Value = 22
$ws.Range("I4").Value = 23
$ws.Range("K4").Value = 91
$ws.Range("L4").Value = 30
$ws.Range("N4").Value = 14
$ws.Range("T4").Value = 177
$ws.Range("U4").Value = 29
$ws.Range("X4").Value = 37
$ws.Range("AA4").Value = 197
$ws.Range("AB4").Value = 15

# Row 5: Dwell time (ms)
$ws.Range("C5").Value = 8600.799999999999
$ws.Range("D5").Value = 96356.44
$ws.Range("E5").Value = 7757.81
$ws.Range("I5").Value = 13081.08
$ws.Range("K5").Value = 35769.38
$ws.Range("L5").Value = 16009.19
$ws.Range("N5").Value = 10660.89
$ws.Range("T5").Value = 69587.72
$ws.Range("U5").Value = 12788.75
$ws.Range("X5").Value = 15917.11
$ws.Range("AA5").Value = 76892.25999999999
$ws.Range("AB5").Value = 7099.77

# Row 6: Dwell time (%)
$ws.Range("C6").Value = 2.74
$ws.Range("D6").Value = 30.68
$ws.Range("E6").Value = 2.48
$ws.Range("F6").Value = 0.99
$ws.Range("G6").Value = 1.65
$ws.Range("H6").Value = 2.2
$ws.Range("I6").Value = 4.17
$ws.Range("J6").Value = 0.59
$ws.Range("K6").Value = 11.41
$ws.Range("L6").Value = 5.11
$ws.Range("M6").Value = 0.18
$ws.Range("N6").Value = 3.4
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 1.04
$ws.Range("T6").Value = 22.2
$ws.Range("U6").Value = 4.08
$ws.Range("V6").Value = 0.25
$ws.Range("X6").Value = 5.08
$ws.Range("Y6").Value = 0.9399999999999999
$ws.Range("Z6").Value = 1.19
$ws.Range("AA6").Value = 24.49
$ws.Range("AB6").Value = 2.27
$ws.Range("AC6").Value = 0.1
$ws.Range("AD6").Value = 0.47
$ws.Range("AE6").Value = 3.03
$ws.Range("AG6").Value = 0.18

# Row 7: Fixation duration (ms)
$ws.Range("C7").Value = 430.04
$ws.Range("D7").Value = 369.18
$ws.Range("E7").Value = 352.63
$ws.Range("I7").Value = 568.74
$ws.Range("K7").Value = 393.07
$ws.Range("L7").Value = 533.64
$ws.Range("N7").Value = 761.49
$ws.Range("T7").Value = 393.15
$ws.Range("U7").Value = 440.99
$ws.Range("X7").Value = 430.19
$ws.Range("AA7").Value = 390.32
$ws.Range("AB7").Value = 473.32

# --- Remove trailing blank row 10 (dimension shrinks to A1:AG9) ---
$ws.Rows.Item(10).Delete()
